$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D4").Value = "Reimbursement, airtime - Cell usage, 2018 April"

$ws.Range("D4").Select()
